# Update cell values in columns A and B (rows 1-32)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(-0.24056878790175773,-0.21827775150838491,-0.11046363447599639,-0.10213155204484714,-0.098705577786568632,0.0019502865433729255,0.012039907084569812,0.0051481731233558037,0.0071678271157114004,0.0091790839983367789,0.012179059432782857,0.015685115840273323,0.019204913930559897,0.027222143397023224,0.028254163363796714,0.030296672896784838,-0.0040034316069528231,-0.016106195586321093,-0.012091858913602316,-0.0080171232407764137,-0.0040057179633441464,-0.10305205126171124,-0.040497303515878436,-0.020098761786426778,-0.097212139425380428,-0.094588834283362289,-0.091930569451236455,-0.088996821143751959,-0.081357430546493248,-0.02117178432631972,-0.014024317199544001,-0.0040015990799009415)
$bValues = @(0.24038394707937982,0.21751568572317215,0.11013155201914593,0.10170557777145817,0.097261852865691445,-0.0020399071202010788,-0.012053265720339823,-0.0051678271314941071,-0.0071790840138503675,-0.0091790594506768741,-0.012180062395697355,-0.01570491394871576,-0.019222143425674965,-0.027254163374895057,-0.02829667290946869,-0.030475629199908205,0.0039999999839865907,0.016091858903095613,0.012017123229582261,0.0080057179520238719,0.0039999999886077831,0.1024291983203014,0.040098761730025245,0.01999999994289503,0.097088834265143475,0.094430569432047307,0.090996821121711591,0.088357430509137025,0.08117178415905979,0.021024317159483719,0.014001599032416934,0.0039999999669504405)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Widen column B to match column A's width
$ws.Columns.Item(2).ColumnWidth = 15.42578125
